$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per latest crypto data refresh.
# D-column values are forced to remain plain text (NumberFormat "@") so that
# values like "65.931.67" or "0.0000112" are not reinterpreted as numbers/
# scientific notation by Excel; the style is reset back to "Normal" afterward
# so no unintended style/formatting change is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.931.67'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.76%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.684.24'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.89%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '605.06'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.57'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.20%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("E8").Value = '  -0.72%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.124'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.07%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.97'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.50%  '
$ws.Range("E11").Value = '  -3.44%  '
$ws.Range("E12").Value = '  +0.53%  '
$ws.Range("E13").Value = '  +6.97%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '29.84'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.94%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.165.97'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.86%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.733.33'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.77%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.686.40'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.27%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.69'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.80%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.85'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.66'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.93%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '355.12'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.94%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.63'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.60%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000112'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +11.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.85'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.14%  '
$ws.Range("E26").Value = '  -6.03%  '
$ws.Range("E27").Value = '  +3.78%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.64'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.81%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.24'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.30%  '
$ws.Range("E30").Value = '  -3.27%  '
$ws.Range("E31").Value = '  -0.07%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '527.86'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.71%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.78'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.52%  '
$ws.Range("E34").Value = '  +2.67%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.44'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.63%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.431'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.88%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '20.61'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.30%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '161.57'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.72%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.98'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.51%  '
$ws.Range("E40").Value = '  +0.05%  '
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '42.49'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.95%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '166.32'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.14%  '
$ws.Range("E44").Value = '  -2.00%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0631'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.48%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '23.51'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.50%  '
$ws.Range("E47").Value = '  -3.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0264'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.39%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.654'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.51%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '20.56'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.29%  '
$ws.Range("E51").Value = '  +0.52%  '
